# "Site Ready For Kevin's mailing tomorrow"
# Mark a batch of players as Attending (column E) ahead of the mailing,
# and fill in a few missing Year (C) / Coach (D) values that came in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark players as Attending (E column: FALSE -> TRUE) ---
$attendingRows = @(23, 70, 85, 88, 119, 121, 143, 166, 167, 178, 182, 183, 189)
foreach ($r in $attendingRows) {
    $ws.Cells.Item($r, 5).Value = $true
}

# --- Fill in newly-known Year / Coach details ---

# Gorski, Gregory (row 88): Year 1971, Coach Robert Casciola
$ws.Cells.Item(88, 3).Value = 1971
$ws.Cells.Item(88, 4).Value = "Robert Casciola"

# Leahy, Bill (row 121): Year 1978
$ws.Cells.Item(121, 3).Value = 1978

# Olzacki, James (row 166): Year 1970
$ws.Cells.Item(166, 3).Value = 1970

# O'Roark, Mike (row 167): Year 1975
$ws.Cells.Item(167, 3).Value = 1975

# Pavasaris, Cris (row 178): Year was blank text, now a real year 2007
$ws.Cells.Item(178, 3).Value = 2007

# Reed, Hank (row 189): Year correction 1975 -> 1976
$ws.Cells.Item(189, 3).Value = 1976

# --- Update the saved view state to where Kevin left off reviewing ---
$ws.Range("Q166").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 139
$excel.ActiveWindow.ScrollColumn = 1
